$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.447.03'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.796.76'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.97%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '338.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.96%  '
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  +1.54%  '
$ws.Range('E8').Value = '  +1.88%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.72'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.203'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07527'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.09'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.479'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.796.80'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.086'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('E17').Value = '  +2.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06670'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '85.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.63%  '
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.544'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.15%  '
$ws.Range('E22').Value = '  +4.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.437.88'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.55'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.421'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.569'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.39%  '
$ws.Range('E27').Value = '  +0.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '152.89'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.000.67'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '133.73'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.047'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.118'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08714'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.94%  '
$ws.Range('E36').Value = '  -2.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.468'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6928'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.80%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06402'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.89%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.897'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2208'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.02349'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.274'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.53'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6469'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.80%  '
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.876'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.139'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.97%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '130.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07201'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.81'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.12%  '
